$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "0.324s"
$ws.Range("E3").Value = "0.311s"
$ws.Range("E4").Value = "0.321s"
$ws.Range("E5").Value = "0.311s"
$ws.Range("E6").Value = "0.322s"
$ws.Range("E7").Value = "0.232s"
$ws.Range("E8").Value = "0.240s"
$ws.Range("E9").Value = "0.231s"
$ws.Range("E10").Value = "0.234s"
$ws.Range("E11").Value = "0.240s"
$ws.Range("E12").Value = "0.456s"
$ws.Range("E13").Value = "0.457s"
$ws.Range("E14").Value = "0.454s"
$ws.Range("E15").Value = "0.229s"
$ws.Range("E16").Value = "0.229s"
$ws.Range("E17").Value = "0.228s"
